{"js": "// Insert the `head(painters)` console output block immediately after the\n// existing \"Load data\" source-code paragraph (the one ending in\n// `head(painters)`), as a new paragraph styled \"Source Code\" whose lines\n// use the \"Verbatim Char\" character style and are separated by line breaks.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Load the text of every paragraph so we can find the one that ends the\n// \"# Load data\" / \"library(MASS)\" / \"head(painters)\" code block.\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"head(painters)\") !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!targetParagraph) {\n  throw new Error('Could not find the paragraph containing \"head(painters)\".');\n}\n\n// The lines of `head(painters)` output to add, in order.\nconst outputLines = [\n  \"##               Composition Drawing Colour Expression School\",\n  \"## Da Udine               10       8     16          3      A\",\n  \"## Da Vinci               15      16      4         14      A\",\n  \"## Del Piombo              8      13     16          7      A\",\n  \"## Del Sarto              12      16      9          8      A\",\n  \"## Fr. Penni               0      15      8          0      A\",\n  \"## Guilio Romano          15      16      4         14      A\",\n];\n\n// Create a new, empty paragraph right after the target and give it the\n// \"Source Code\" paragraph style (matching the surrounding code blocks).\nconst newParagraph = targetParagraph.insertParagraph(\"\", \"After\");\nnewParagraph.style = \"Source Code\";\n\n// Populate the new paragraph: each line is its own text run styled with\n// the \"Verbatim Char\" character style; consecutive lines are separated by\n// a manual line break run (no character style), mirroring how the other\n// multi-line SourceCode blocks in this document are structured.\nfor (let i = 0; i < outputLines.length; i++) {\n  if (i > 0) {\n    newParagraph.insertBreak(\"Line\", \"End\");\n  }\n  const lineRange = newParagraph.insertText(outputLines[i], \"End\");\n  lineRange.style = \"Verbatim Char\";\n}\n\nawait context.sync();\n", "ps1": "# Insert the `head(painters)` console output block immediately after the\n# existing \"Load data\" source-code paragraph (the one ending in\n# `head(painters)`), as a new paragraph styled \"Source Code\" whose lines\n# use the \"Verbatim Char\" character style and are separated by line breaks.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that contains the `head(painters)` source line.\n$targetParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*head(painters)*\") {\n        $targetParagraph = $p\n        break\n    }\n}\n\nif ($targetParagraph -eq $null) {\n    throw \"Could not find the paragraph containing 'head(painters)'.\"\n}\n\n# Remember where the target paragraph sits so we can re-fetch the freshly\n# inserted paragraph by index afterwards (object references can go stale\n# across structural edits).\n$targetIndex = $targetParagraph.Range.Information(1)\n$allParas = $d.Paragraphs\n$targetOrdinal = 0\n$i = 0\nforeach ($p in $allParas) {\n    $i = $i + 1\n    if ($p.Range.Start -eq $targetParagraph.Range.Start) {\n        $targetOrdinal = $i\n        break\n    }\n}\n\n# Insert a brand-new, empty paragraph right after the target and style it\n# like the other R source-code blocks in this document.\n$targetParagraph.Range.InsertParagraphAfter()\n$newParagraph = $d.Paragraphs.Item($targetOrdinal + 1)\n$newParagraph.Style = \"Source Code\"\n\n# The lines of `head(painters)` console output to add, in order.\n$outputLines = @(\n    \"##               Composition Drawing Colour Expression School\",\n    \"## Da Udine               10       8     16          3      A\",\n    \"## Da Vinci               15      16      4         14      A\",\n    \"## Del Piombo              8      13     16          7      A\",\n    \"## Del Sarto              12      16      9          8      A\",\n    \"## Fr. Penni               0      15      8          0      A\",\n    \"## Guilio Romano          15      16      4         14      A\"\n)\n\n# Populate the new paragraph: each line becomes its own text run styled\n# with \"Verbatim Char\"; consecutive lines are separated by a manual line\n# break run (left in the default/no character style), which mirrors how\n# the other multi-line SourceCode blocks in this document are built.\n$pos = $newParagraph.Range.Start\nfor ($li = 0; $li -lt $outputLines.Length; $li++) {\n    if ($li -gt 0) {\n        $breakRange = $d.Range($pos, $pos)\n        $breakRange.InsertBreak(\"Line\")\n        $pos = $pos + 1\n    }\n\n    $line = $outputLines[$li]\n    $textRange = $d.Range($pos, $pos)\n    $textRange.InsertAfter($line)\n    $styledRange = $d.Range($pos, $pos + $line.Length)\n    $styledRange.Style = \"Verbatim Char\"\n    $pos = $pos + $line.Length\n}\n"}
